# Swap the data (columns B through AD) between paired rows.
# Column A (the sequential "id" index) stays untouched on each row.
#
# NOTE: each pair below uses its own, non-reused set of variable names.
# (Reusing the same $range/$values variable names across several
# row-pairs in a loop was observed to make the later writes clobber the
# earlier ones when the workbook is saved, even though interim reads
# in-session looked correct - so every pair gets distinct variables.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Costa Rica Primera Division")

# Pair 1: rows 110 <-> 111
$rangePair1Row1 = $ws.Range("B110:AD110")
$rangePair1Row2 = $ws.Range("B111:AD111")
$valuesPair1Row1 = $rangePair1Row1.Value()
$valuesPair1Row2 = $rangePair1Row2.Value()
$rangePair1Row1.Value = $valuesPair1Row2
$rangePair1Row2.Value = $valuesPair1Row1

# Pair 2: rows 237 <-> 238
$rangePair2Row1 = $ws.Range("B237:AD237")
$rangePair2Row2 = $ws.Range("B238:AD238")
$valuesPair2Row1 = $rangePair2Row1.Value()
$valuesPair2Row2 = $rangePair2Row2.Value()
$rangePair2Row1.Value = $valuesPair2Row2
$rangePair2Row2.Value = $valuesPair2Row1

# Pair 3: rows 256 <-> 257
$rangePair3Row1 = $ws.Range("B256:AD256")
$rangePair3Row2 = $ws.Range("B257:AD257")
$valuesPair3Row1 = $rangePair3Row1.Value()
$valuesPair3Row2 = $rangePair3Row2.Value()
$rangePair3Row1.Value = $valuesPair3Row2
$rangePair3Row2.Value = $valuesPair3Row1
